
# shared_count_record.xlsx: "updated some datasheet terms for clarity"
# The "Plant" column (V) values are renamed:
#   Jasmine          -> M. paniculata
#   Curry            -> B. koenigii
#   Mix (OJ/Curry)   -> Mix

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yield")

$plantCol = $ws.Range("V1:V71")

# Use whole-cell matching (xlWhole = 1) so "Mix (OJ/Curry)" is not
# partially clobbered by the "Curry" replacement.
$plantCol.Replace("Jasmine", "M. paniculata", 1) | Out-Null
$plantCol.Replace("Curry", "B. koenigii", 1) | Out-Null
$plantCol.Replace("Mix (OJ/Curry)", "Mix", 1) | Out-Null

# Restore the active cell/selection recorded for this sheet.
$ws.Activate()
$ws.Range("V29").Select()
